$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 519.93335
$ws.Range("J58").Value = 1766.6666
$ws.Range("L58").Value = 5299.9998
$ws.Range("N58").Value = -5599.9998

$ws.Range("H61").Value = 1015.3
$ws.Range("I61").Value = 165
$ws.Range("J61").Value = 2999.3333
$ws.Range("K61").Value = 495
$ws.Range("L61").Value = 8997.999899999999
$ws.Range("M61").Value = -323
$ws.Range("N61").Value = -9341.999899999999

$ws.Range("H87").Value = 31257.143
$ws.Range("J87").Value = 31257.143
$ws.Range("L87").Value = 31257.143
$ws.Range("N87").Value = -33753.143

$ws.Range("H90").Value = 31257.143
$ws.Range("J90").Value = 31257.143
$ws.Range("L90").Value = 93771.429
$ws.Range("N90").Value = -106251.429

$ws.Range("H129").Value = 549.9545000000001
$ws.Range("J129").Value = 1044.4286
$ws.Range("L129").Value = 3133.2858
$ws.Range("N129").Value = -13133.2858

$ws.Range("H138").Value = 2643.328
$ws.Range("I138").Value = 1337.0264
$ws.Range("J138").Value = 4801.5654
$ws.Range("K138").Value = 4011.0792
$ws.Range("L138").Value = 14404.6962
$ws.Range("M138").Value = 1128.9208
$ws.Range("N138").Value = -24684.6962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3100.65
$ws.Range("I32").Value = 2135.247
$ws.Range("J32").Value = 10911.637
$ws.Range("K32").Value = 2135.247
$ws.Range("L32").Value = 10911.637
$ws.Range("M32").Value = -1848.247
$ws.Range("N32").Value = -11485.637

$ws.Range("H63").Value = 2995.6667
$ws.Range("I63").Value = 3450.7144
$ws.Range("K63").Value = 3450.7144
$ws.Range("M63").Value = -2764.7144

$ws.Range("H66").Value = 2995.6667
$ws.Range("I66").Value = 3450.7144
$ws.Range("K66").Value = 17253.572
$ws.Range("M66").Value = -13821.572

$ws.Range("H122").Value = 2759.4
$ws.Range("I122").Value = 2294
$ws.Range("J122").Value = 3457.5
$ws.Range("K122").Value = 6882
$ws.Range("L122").Value = 10372.5
$ws.Range("M122").Value = -4432
$ws.Range("N122").Value = -15272.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1458.3334
$ws.Range("I134").Value = 1182.2
$ws.Range("J134").Value = 1803.5
$ws.Range("K134").Value = 3546.6
$ws.Range("L134").Value = 5410.5
$ws.Range("M134").Value = -1011.6
$ws.Range("N134").Value = -10480.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 703.625
$ws.Range("I16").Value = 662.25
$ws.Range("J16").Value = 745
$ws.Range("K16").Value = 662.25
$ws.Range("L16").Value = 745
$ws.Range("M16").Value = -375.25
$ws.Range("N16").Value = -1319

$ws.Range("H31").Value = 31807.057
$ws.Range("I31").Value = 48770.363
$ws.Range("J31").Value = 3099.923
$ws.Range("K31").Value = 48770.363
$ws.Range("L31").Value = 3099.923
$ws.Range("M31").Value = -48475.363
$ws.Range("N31").Value = -3689.923

$ws.Range("H34").Value = 31807.057
$ws.Range("I34").Value = 48770.363
$ws.Range("J34").Value = 3099.923
$ws.Range("K34").Value = 48770.363
$ws.Range("L34").Value = 3099.923
$ws.Range("M34").Value = -48568.363
$ws.Range("N34").Value = -3503.923

$ws.Range("H58").Value = 1008.32355
$ws.Range("I58").Value = 888.902
$ws.Range("J58").Value = 1366.5883
$ws.Range("K58").Value = 888.902
$ws.Range("L58").Value = 1366.5883
$ws.Range("M58").Value = -685.902
$ws.Range("N58").Value = -1772.5883

$ws.Range("H113").Value = 703.625
$ws.Range("I113").Value = 662.25
$ws.Range("J113").Value = 745
$ws.Range("K113").Value = 662.25
$ws.Range("L113").Value = 745
$ws.Range("M113").Value = 1507.75
$ws.Range("N113").Value = -5085

$ws.Range("H122").Value = 2182.7585
$ws.Range("I122").Value = 3112.6667
$ws.Range("J122").Value = 1186.4286
$ws.Range("K122").Value = 9338.000100000001
$ws.Range("L122").Value = 3559.2858
$ws.Range("M122").Value = -6888.000100000001
$ws.Range("N122").Value = -8459.2858

$ws.Range("H136").Value = 1008.32355
$ws.Range("I136").Value = 888.902
$ws.Range("J136").Value = 1366.5883
$ws.Range("K136").Value = 2666.706
$ws.Range("L136").Value = 4099.7649
$ws.Range("M136").Value = -116.7060000000001
$ws.Range("N136").Value = -9199.7649

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1495.3939
$ws.Range("I5").Value = 357.23077
$ws.Range("J5").Value = 2235.2
$ws.Range("K5").Value = 1071.69231
$ws.Range("L5").Value = 6705.599999999999
$ws.Range("M5").Value = -959.6923099999999
$ws.Range("N5").Value = -6929.599999999999

$ws.Range("H107").Value = 354.95456
$ws.Range("I107").Value = 177.42857
$ws.Range("J107").Value = 437.8
$ws.Range("K107").Value = 532.28571
$ws.Range("L107").Value = 1313.4
$ws.Range("M107").Value = 1387.71429
$ws.Range("N107").Value = -5153.4

$ws.Range("H122").Value = 1132.2632
$ws.Range("I122").Value = 846.86957
$ws.Range("J122").Value = 1569.8667
$ws.Range("K122").Value = 7621.826129999999
$ws.Range("L122").Value = 14128.8003
$ws.Range("M122").Value = -5171.826129999999
$ws.Range("N122").Value = -19028.8003

$ws.Range("H135").Value = 1495.3939
$ws.Range("I135").Value = 357.23077
$ws.Range("J135").Value = 2235.2
$ws.Range("K135").Value = 3215.07693
$ws.Range("L135").Value = 20116.8
$ws.Range("M135").Value = -680.0769300000002
$ws.Range("N135").Value = -25186.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 80000
$ws.Range("J26").Value = 80000
$ws.Range("L26").Value = 80000
$ws.Range("N26").Value = -80560

$ws.Range("H50").Value = 80000
$ws.Range("J50").Value = 80000
$ws.Range("L50").Value = 80000
$ws.Range("N50").Value = -80996

$ws.Range("H102").Value = 1356.5333
$ws.Range("I102").Value = 1385.5714
$ws.Range("J102").Value = 950
$ws.Range("K102").Value = 1385.5714
$ws.Range("L102").Value = 950
$ws.Range("M102").Value = 236.4286
$ws.Range("N102").Value = -4194

$ws.Range("H107").Value = 606.7143
$ws.Range("I107").Value = 483.83334
$ws.Range("K107").Value = 483.83334
$ws.Range("M107").Value = 1436.16666

$ws.Range("H122").Value = 1675.129
$ws.Range("I122").Value = 1288.0435
$ws.Range("J122").Value = 2788
$ws.Range("K122").Value = 3864.1305
$ws.Range("L122").Value = 8364
$ws.Range("M122").Value = -1414.1305
$ws.Range("N122").Value = -13264

$ws.Range("H126").Value = 1398.4286
$ws.Range("I126").Value = 1256
$ws.Range("J126").Value = 1455.4
$ws.Range("K126").Value = 3768
$ws.Range("L126").Value = 4366.200000000001
$ws.Range("M126").Value = -1298
$ws.Range("N126").Value = -9306.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1631.875
$ws.Range("I7").Value = 1161.25
$ws.Range("J7").Value = 2102.5
$ws.Range("K7").Value = 1161.25
$ws.Range("L7").Value = 2102.5
$ws.Range("M7").Value = -1049.25
$ws.Range("N7").Value = -2326.5

$ws.Range("H40").Value = 3551.9412
$ws.Range("I40").Value = 3237.5
$ws.Range("J40").Value = 5019.3335
$ws.Range("K40").Value = 3237.5
$ws.Range("L40").Value = 5019.3335
$ws.Range("M40").Value = -3101.5
$ws.Range("N40").Value = -5291.3335

$ws.Range("H46").Value = 975525.25
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 3900001
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 3900001
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -3900377

$ws.Range("H122").Value = 3191.5186
$ws.Range("I122").Value = 3315.8462
$ws.Range("J122").Value = 3076.0715
$ws.Range("K122").Value = 9947.5386
$ws.Range("L122").Value = 9228.2145
$ws.Range("M122").Value = -7497.5386
$ws.Range("N122").Value = -14128.2145

$ws.Range("H126").Value = 1631.875
$ws.Range("I126").Value = 1161.25
$ws.Range("J126").Value = 2102.5
$ws.Range("K126").Value = 3483.75
$ws.Range("L126").Value = 6307.5
$ws.Range("M126").Value = -1013.75
$ws.Range("N126").Value = -11247.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 418323.34
$ws.Range("I122").Value = 527734.7
$ws.Range("J122").Value = 2560.2
$ws.Range("K122").Value = 1583204.1
$ws.Range("L122").Value = 7680.599999999999
$ws.Range("M122").Value = -1580754.1
$ws.Range("N122").Value = -12580.6

$ws.Range("H126").Value = 400794.8
$ws.Range("I126").Value = 500661.66
$ws.Range("J126").Value = 1327.4
$ws.Range("K126").Value = 1501984.98
$ws.Range("L126").Value = 3982.2
$ws.Range("M126").Value = -1499514.98
$ws.Range("N126").Value = -8922.200000000001
